$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data rows appended to the bottom of the table (rows 80-82)
# Column A holds date-formatted serial values (style already applied via AutoFill from row above)

$newRows = @(
    @{ Row = 80; DateSerial = 44443; Values = @(121.1, 228.1, 264.5, 274.2, 218.6, 223.1, 192, 142.9, 109.3, 99.1) },
    @{ Row = 81; DateSerial = 44450; Values = @(89.6, 165.7, 192.3, 198.8, 162.1, 166.1, 145.80000000000001, 109.1, 82.8, 75.5) },
    @{ Row = 82; DateSerial = 44457; Values = @(7.7, 13.5, 17.7, 18.2, 17.600000000000001, 15.8, 13.8, 10.4, 8.3000000000000007, 7.3) }
)

# Carry the date style (column A) down from the last existing row onto the
# new rows so the added date cells keep the same date number-formatting.
$ws.Range("A79").Copy()
$ws.Range("A80:A82").PasteSpecial(-4122) # xlPasteFormats
$ws.Application.CutCopyMode = 0

foreach ($entry in $newRows) {
    $r = $entry.Row

    # Column A: date value (style was already copied above).
    $ws.Cells.Item($r, 1).Value = $entry.DateSerial

    # Columns B..K: numeric values
    for ($i = 0; $i -lt $entry.Values.Count; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($r, $col).Value = $entry.Values[$i]
    }
}

# Update the view: scroll down and move the active selection to reflect
# where the user ended up after entering the new data.
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("F95").Select()
